$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bText = @"
Pipeline(steps=[('scaler',
                 ColumnTransformer(n_jobs=-1, remainder='passthrough',
                                   transformers=[('StandardScaler',
                                                  StandardScaler(),
                                                  ['AE_HR', 'AE_V',
                                                   'AbsOffAxis_HR',
                                                   'AbsOffAxis_V',
                                                   'AbsOnAxis_HR',
                                                   'AbsOnAxis_V', 'BallPath_HR',
                                                   'BallPath_V', 'CMT_HR',
                                                   'CMT_V', 'Corrective_HR',
                                                   'Corrective_V', 'Delta_AE',
                                                   'Delta_Fullpath', 'Delta_MT',
                                                   'Delta_OffAxis',
                                                   'Delta_OnAxis', 'Delta_PV',
                                                   'Delta_RT', 'FullPath_HR',
                                                   'FullPath_V', 'MT_HR',
                                                   'MT_V', 'PeakV_HR',
                                                   'PeakV_V', 'RT_HR', 'RT_V',
                                                   'TMT_HR', 'TMT_V', 'VE_HR', ...])])),
                ('selector', RandomUnderSampler(random_state=42)),
                ('model',
                 DecisionTreeClassifier(class_weight='balanced', max_depth=4,
                                        max_features='sqrt',
                                        min_samples_leaf=11,
                                        min_samples_split=13,
                                        random_state=42))])
"@

$cText = @"
{'selector': RandomUnderSampler(random_state=42), 'scaler': ColumnTransformer(n_jobs=-1, remainder='passthrough',
                  transformers=[('StandardScaler', StandardScaler(),
                                 ['AE_HR', 'AE_V', 'AbsOffAxis_HR',
                                  'AbsOffAxis_V', 'AbsOnAxis_HR', 'AbsOnAxis_V',
                                  'BallPath_HR', 'BallPath_V', 'CMT_HR',
                                  'CMT_V', 'Corrective_HR', 'Corrective_V',
                                  'Delta_AE', 'Delta_Fullpath', 'Delta_MT',
                                  'Delta_OffAxis', 'Delta_OnAxis', 'Delta_PV',
                                  'Delta_RT', 'FullPath_HR', 'FullPath_V',
                                  'MT_HR', 'MT_V', 'PeakV_HR', 'PeakV_V',
                                  'RT_HR', 'RT_V', 'TMT_HR', 'TMT_V', 'VE_HR', ...])]), 'model__min_samples_split': 13, 'model__min_samples_leaf': 11, 'model__max_features': 'sqrt', 'model__max_depth': 4, 'model__criterion': 'gini', 'model__class_weight': 'balanced'}
"@

$ws.Range("B2").Value = $bText
$ws.Range("C2").Value = $cText
$ws.Range("D2").Value = 0.5387144812980054
$ws.Range("F2").Value = 0.7121122738459804
$ws.Range("G2").Value = 0.5833333333333334
$ws.Range("I2").Value = "[0 1 1 1 1 0 1 1 0 1 1 1 1 0 1 0 0 0 1 0 0 0 1 1]"

$ws.Range("B3").Value = $bText
$ws.Range("C3").Value = $cText
$ws.Range("D3").Value = 0.5236397472795373
$ws.Range("F3").Value = 0.7275612534486958
$ws.Range("G3").Value = 0.6967109424414927
$ws.Range("I3").Value = "[1 1 1 0 1 1 1 0 0 1 0 0 1 0 1 0 1 1 1 1 1 1 1 1]"

$ws.Range("B4").Value = $bText
$ws.Range("C4").Value = $cText
$ws.Range("D4").Value = 0.4535291717077287
$ws.Range("F4").Value = 0.6819285466527922
$ws.Range("G4").Value = 0.6458333333333334
$ws.Range("I4").Value = "[0 1 1 1 0 1 1 1 1 1 1 1 0 1 1 1 0 0 1 1 1 1 0 1]"

$ws.Range("B5").Value = $bText
$ws.Range("C5").Value = $cText
$ws.Range("D5").Value = 0.5797953697937545
$ws.Range("F5").Value = 0.7435179549604384
$ws.Range("G5").Value = 0.4207459207459207
$ws.Range("I5").Value = "[0 0 1 1 1 0 0 0 0 1 0 1 1 0 0 1 1 1 1 0 1 0 0 1]"

$ws.Range("B6").Value = $bText
$ws.Range("C6").Value = $cText
$ws.Range("D6").Value = 0.5230189578658997
$ws.Range("F6").Value = 0.7013753438359589
$ws.Range("G6").Value = 0.6666666666666666
$ws.Range("I6").Value = "[0 0 1 1 1 0 0 0 1 1 1 1 0 1 1 1 0 0 0 1 0 1 1 1]"

